{"js": "// 1) Title text: \"Titulo do documento (Principal)\" -> \"Associa\u00e7\u00e3o de ca\u00e7a\"\nconst titleResults = context.document.body.search(\"Titulo do documento (Principal)\", { matchCase: true });\ntitleResults.load(\"text\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"Associa\u00e7\u00e3o de ca\u00e7a\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) \"MathTy\" + \"pe\" (two runs) -> single run \"MathType\"\nconst mathTypeResults = context.document.body.search(\"MathType\", { matchCase: true });\nmathTypeResults.load(\"text\");\nawait context.sync();\nif (mathTypeResults.items.length > 0) {\n  mathTypeResults.items[0].insertText(\"MathType\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) \"M_T\" + \"ext\" (two runs) -> single run \"M_Text\".\n// There are two paragraphs in the document that read \"Main text paragraph (M_Text).\"\n// right after the \"3. Experimental Section (M_Heading1)\" heading; the second one\n// is the one whose \"M_Text\" is still split across two runs. Anchor on the unique\n// heading text and walk forward two paragraphs to reach it precisely.\nconst headingResults = context.document.body.search(\"3. Experimental Section\", { matchCase: true });\nawait context.sync();\nif (headingResults.items.length > 0) {\n  const headingParagraphs = headingResults.items[0].paragraphs;\n  headingParagraphs.load(\"text\");\n  await context.sync();\n\n  const headingParagraph = headingParagraphs.items[0];\n  const firstFollowingParagraph = headingParagraph.next();\n  const secondFollowingParagraph = firstFollowingParagraph.next();\n  const targetRange = secondFollowingParagraph.getRange();\n\n  const mTextResults = targetRange.search(\"M_Text\", { matchCase: true });\n  mTextResults.load(\"text\");\n  await context.sync();\n  if (mTextResults.items.length > 0) {\n    mTextResults.items[0].insertText(\"M_Text\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 4) \", Country, 2008\" + \"; \" (two runs) -> single run \", Country, 2008; \"\nconst countryResults = context.document.body.search(\", Country, 2008; \", { matchCase: true });\ncountryResults.load(\"text\");\nawait context.sync();\nif (countryResults.items.length > 0) {\n  countryResults.items[0].insertText(\", Country, 2008; \", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Title text: \"Titulo do documento (Principal)\" -> \"Associa\u00e7\u00e3o de ca\u00e7a\"\n$titleFind = $d.Content.Find\n$titleFind.Text = \"Titulo do documento (Principal)\"\n$titleFind.Replacement.Text = \"Associa\u00e7\u00e3o de ca\u00e7a\"\n$titleFind.Execute($titleFind.Text, $false, $false, $false, $false, $false, $true, 0, $false, $titleFind.Replacement.Text, 1) | Out-Null\n\n# 2) \"MathTy\" + \"pe\" (two runs) -> single run \"MathType\".\n# Find.Execute matches across the run boundary and rewrites it as one run.\n$mathTypeFind = $d.Content.Find\n$mathTypeFind.Text = \"MathType\"\n$mathTypeFind.Replacement.Text = \"MathType\"\n$mathTypeFind.Execute($mathTypeFind.Text, $false, $false, $false, $false, $false, $true, 0, $false, $mathTypeFind.Replacement.Text, 1) | Out-Null\n\n# 3) \"M_T\" + \"ext\" (two runs) -> single run \"M_Text\".\n# The document has two paragraphs reading \"Main text paragraph (M_Text).\" right\n# after the \"3. Experimental Section (M_Heading1)\" heading; only the second one\n# still has \"M_Text\" split across two runs. Anchor on the unique heading text,\n# then scan \"M_Text\" matches that occur after it and stop on the 2nd one.\n$headingRange = $d.Content\n$headingFind = $headingRange.Find\n$headingFind.Text = \"3. Experimental Section\"\n$headingFind.Forward = $true\n$headingFind.Wrap = 0\n$headingFind.Execute() | Out-Null\n$anchorStart = $headingRange.Start\n\n$scanRange = $d.Content\n$scanFind = $scanRange.Find\n$matchesAfterHeading = 0\n$targetFound = $false\nwhile (-not $targetFound) {\n    $scanFind.Text = \"M_Text\"\n    $scanFind.Forward = $true\n    $scanFind.Wrap = 0\n    $found = $scanFind.Execute()\n    if (-not $found) { break }\n    if ($scanRange.Start -gt $anchorStart) {\n        $matchesAfterHeading = $matchesAfterHeading + 1\n        if ($matchesAfterHeading -eq 2) {\n            # Re-assigning the exact same text is a no-op for the underlying\n            # run merge, so stage a different value first, then set the real\n            # text; $scanRange automatically re-tracks its own bounds.\n            $scanRange.Text = \"M_Text__TEMP__\"\n            $scanRange.Text = \"M_Text\"\n            $targetFound = $true\n        }\n    }\n}\n\n# 4) \", Country, 2008\" + \"; \" (two runs) -> single run \", Country, 2008; \"\n$countryFind = $d.Content.Find\n$countryFind.Text = \", Country, 2008; \"\n$countryFind.Replacement.Text = \", Country, 2008; \"\n$countryFind.Execute($countryFind.Text, $false, $false, $false, $false, $false, $true, 0, $false, $countryFind.Replacement.Text, 1) | Out-Null\n"}
